# Commit: Update standby badges to be yellow for better visibility
# Data-side change that accompanies the badge re-color: Peter Adamidis's
# confirmed seat assignment (row 2 of "Seat Assignments") is converted into
# a pending "Standby" record, which now lives on its own new "Standbys"
# sheet (inserted right before "Groups").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "Standbys" worksheet right before "Groups".
# ---------------------------------------------------------------------------
$groupsSheet = $wb.Worksheets.Item("Groups")
$standbys = $wb.Worksheets.Add($groupsSheet)
$standbys.Name = "Standbys"

$standbys.Range("A1").Value = "ID"
$standbys.Range("B1").Value = "RecordDayID"
$standbys.Range("C1").Value = "ContestantID"
$standbys.Range("D1").Value = "Status"
$standbys.Range("E1").Value = "Notes"

$standbys.Range("A2").Value = "69511596-e3b5-41be-a93b-920748af4fe0"
$standbys.Range("B2").Value = "e432f0fe-1383-44a2-990c-5f787da5008a"
$standbys.Range("C2").Value = "0ccaf8bc-6ade-4ad6-9537-92f9829a6502"
$standbys.Range("D2").Value = "pending"

# ---------------------------------------------------------------------------
# 2) Contestants sheet: rows rotate - the row order was resorted (by name),
#    so each contestant's row now holds the data that used to sit on the
#    next row (cyclically).
# ---------------------------------------------------------------------------
$contestants = $wb.Worksheets.Item("Contestants")

$contestants.Range("A2").Value = "28603f95-d5f6-47ab-88c4-0d79742a6b02"
$contestants.Range("B2").Value = "Felicity Parker-Hill"
$contestants.Range("C2").Value = 27
$contestants.Range("D2").Value = "Not Specified"
$contestants.Range("E2").Value = "felicity.parkerhill@endemolshine.com.au"
$contestants.Range("F2").Value = "498086080"
$contestants.Range("G2").Value = "Melbourne"
$contestants.Range("I2").Value = "assigned"
$contestants.Range("J2").Value = "Peter Adamidis, Kathleen Reynolds"
$contestants.Range("K2").Value = "5fe641da-4067-49a7-bae7-e63413b3e404"
$contestants.Range("L2").Value = "N"
$contestants.Range("M2").Value = "N/A"

$contestants.Range("A3").Value = "d698b1de-6641-45c6-aa63-f577d2b634bb"
$contestants.Range("B3").Value = "Kathleen Reynolds"
$contestants.Range("C3").Value = 33
$contestants.Range("D3").Value = "Not Specified"
$contestants.Range("E3").Value = "kathleenmonicareynolds@gmail.com"
$contestants.Range("F3").Value = "498086080"
$contestants.Range("G3").Value = "Footscray"
$contestants.Range("I3").Value = "assigned"
$contestants.Range("J3").Value = "Peter Adamidis, Felicity Parker-Hill"
$contestants.Range("K3").Value = "5fe641da-4067-49a7-bae7-e63413b3e404"
$contestants.Range("L3").Value = "N"
$contestants.Range("M3").Value = "N/A"

$contestants.Range("A4").Value = "0ccaf8bc-6ade-4ad6-9537-92f9829a6502"
$contestants.Range("B4").Value = "Peter Adamidis"
$contestants.Range("C4").Value = 34
$contestants.Range("D4").Value = "Not Specified"
$contestants.Range("E4").Value = "peter.adamidis@gmail.com"
$contestants.Range("F4").Value = "498086080"
$contestants.Range("G4").Value = ""
$contestants.Range("I4").Value = "assigned"
$contestants.Range("J4").Value = "Kathleen Reynolds, Felicity Parker-Hill"
$contestants.Range("K4").Value = "5fe641da-4067-49a7-bae7-e63413b3e404"
$contestants.Range("L4").Value = "Y"
$contestants.Range("M4").Value = "Broken Leg"

# ---------------------------------------------------------------------------
# 3) Seat Assignments: Peter Adamidis's seat-assignment row (row 2) is
#    removed - he is now tracked on the Standbys sheet instead. Deleting the
#    row shifts the remaining two rows up and fixes the dimension/used range.
# ---------------------------------------------------------------------------
$seatAssignments = $wb.Worksheets.Item("Seat Assignments")
$seatAssignments.Rows.Item(2).Delete()
